$wb = $excel.ActiveWorkbook

# Rename the "shortDescription" column header to "description" on every sheet (cell C1)
foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "description"
}

# German sheet: move the selection, no longer the active/selected tab
$wsGerman = $wb.Worksheets.Item("German")
$wsGerman.Range("C2").Select()

# English sheet: move the selection and make it the active/selected tab
# (selecting this range last leaves its sheet as the active tab)
$wsEnglish = $wb.Worksheets.Item("English")
$wsEnglish.Range("I11").Select()
